$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1432551
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1432551
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4297653
$ws.Range("N17").Value = -4297989

$ws.Range("H33").Value = 245.75
$ws.Range("I33").Value = 202.88889
$ws.Range("J33").Value = 280.81818
$ws.Range("K33").Value = 202.88889
$ws.Range("L33").Value = 280.81818
$ws.Range("M33").Value = 26.11111
$ws.Range("N33").Value = -738.81818

$ws.Range("H62").Value = 3785.1482
$ws.Range("I62").Value = 2788.2354
$ws.Range("J62").Value = 5479.9
$ws.Range("K62").Value = 2788.2354
$ws.Range("L62").Value = 5479.9
$ws.Range("M62").Value = -2164.2354
$ws.Range("N62").Value = -6727.9

$ws.Range("H65").Value = 3785.1482
$ws.Range("I65").Value = 2788.2354
$ws.Range("J65").Value = 5479.9
$ws.Range("K65").Value = 13941.177
$ws.Range("L65").Value = 27399.5
$ws.Range("M65").Value = -10821.177
$ws.Range("N65").Value = -33639.5

$ws.Range("H98").Value = 765.5714
$ws.Range("I98").Value = 753.35
$ws.Range("J98").Value = 796.125
$ws.Range("K98").Value = 753.35
$ws.Range("L98").Value = 796.125
$ws.Range("M98").Value = 744.65

$ws.Range("H112").Value = 1247.7407
$ws.Range("I112").Value = 879.8
$ws.Range("J112").Value = 1331.3636
$ws.Range("K112").Value = 2639.4
$ws.Range("L112").Value = 3994.0908
$ws.Range("M112").Value = -1531.4
$ws.Range("N112").Value = -6210.0908

$ws.Range("H122").Value = 765.5714
$ws.Range("I122").Value = 753.35
$ws.Range("J122").Value = 796.125
$ws.Range("K122").Value = 2260.05
$ws.Range("L122").Value = 2388.375
$ws.Range("M122").Value = 189.9499999999998

$ws.Range("H129").Value = 1005.3333
$ws.Range("I129").Value = 599.8333
$ws.Range("J129").Value = 1037.7733
$ws.Range("K129").Value = 1799.4999
$ws.Range("L129").Value = 3113.3199
$ws.Range("M129").Value = 3200.5001
$ws.Range("N129").Value = -13113.3199

$ws.Range("H131").Value = 4055.7188
$ws.Range("I131").Value = 934
$ws.Range("J131").Value = 4776.115
$ws.Range("K131").Value = 2802
$ws.Range("L131").Value = 14328.345
$ws.Range("M131").Value = 2238
$ws.Range("N131").Value = -24408.345

$ws.Range("H132").Value = 14373.855
$ws.Range("I132").Value = 16244.137
$ws.Range("J132").Value = 2030
$ws.Range("K132").Value = 48732.411
$ws.Range("L132").Value = 6090
$ws.Range("M132").Value = -46202.411
$ws.Range("N132").Value = -11150

$ws.Range("H138").Value = 4472.659
$ws.Range("I138").Value = 2756.0938
$ws.Range("J138").Value = 5403.6777
$ws.Range("K138").Value = 8268.2814
$ws.Range("L138").Value = 16211.0331
$ws.Range("M138").Value = -3128.2814
$ws.Range("N138").Value = -26491.0331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9494.610000000001
$ws.Range("I32").Value = 7555.5
$ws.Range("J32").Value = 44398.6
$ws.Range("K32").Value = 7555.5
$ws.Range("L32").Value = 44398.6
$ws.Range("M32").Value = -7268.5

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H122").Value = 2481.75
$ws.Range("I122").Value = 1541.1538
$ws.Range("J122").Value = 4228.5713
$ws.Range("K122").Value = 4623.4614
$ws.Range("L122").Value = 12685.7139
$ws.Range("M122").Value = -2173.4614

$ws.Range("H124").Value = 36759.8
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 36759.8
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 36759.8
$ws.Range("N124").Value = -46579.8

$ws.Range("H125").Value = 41821.668
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 41821.668
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 41821.668
$ws.Range("N125").Value = -51661.668

$ws.Range("H132").Value = 1524.4271
$ws.Range("I132").Value = 1164.2715
$ws.Range("J132").Value = 2494.077
$ws.Range("K132").Value = 3492.8145
$ws.Range("L132").Value = 7482.231000000001
$ws.Range("M132").Value = -962.8145000000004
$ws.Range("N132").Value = -12542.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 301.07144
$ws.Range("I22").Value = 202.09091
$ws.Range("J22").Value = 664
$ws.Range("K22").Value = 202.09091
$ws.Range("L22").Value = 664
$ws.Range("M22").Value = -29.09091000000001
$ws.Range("N22").Value = -1010

$ws.Range("H74").Value = 20780
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 20780
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 20780
$ws.Range("N74").Value = -22652

$ws.Range("H77").Value = 20780
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 20780
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 62340
$ws.Range("N77").Value = -71700

$ws.Range("H132").Value = 45526
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 45526
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 45526
$ws.Range("N132").Value = -55646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2634.7
$ws.Range("I31").Value = 1876.7042
$ws.Range("J31").Value = 4490.483
$ws.Range("K31").Value = 1876.7042
$ws.Range("L31").Value = 4490.483
$ws.Range("M31").Value = -1581.7042
$ws.Range("N31").Value = -5080.483

$ws.Range("H34").Value = 2634.7
$ws.Range("I34").Value = 1876.7042
$ws.Range("J34").Value = 4490.483
$ws.Range("K34").Value = 1876.7042
$ws.Range("L34").Value = 4490.483
$ws.Range("M34").Value = -1674.7042
$ws.Range("N34").Value = -4894.483

$ws.Range("H132").Value = 1001932.06
$ws.Range("I132").Value = 1786162.1
$ws.Range("J132").Value = 3821.0908
$ws.Range("K132").Value = 5358486.300000001
$ws.Range("L132").Value = 11463.2724
$ws.Range("M132").Value = -5355956.300000001
$ws.Range("N132").Value = -16523.2724

$ws.Range("H141").Value = 48717.453
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 48717.453
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 48717.453
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -59077.453

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 298.8889
$ws.Range("I4").Value = 241.42857
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 724.28571
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -612.28571
$ws.Range("N4").Value = -1724

$ws.Range("H39").Value = 16706.143
$ws.Range("I39").Value = 50671.5
$ws.Range("J39").Value = 3120
$ws.Range("K39").Value = 152014.5
$ws.Range("L39").Value = 9360
$ws.Range("M39").Value = -151720.5
$ws.Range("N39").Value = -9948

$ws.Range("H46").Value = 266.33334
$ws.Range("I46").Value = 119.6
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 358.8
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -267.8
$ws.Range("N46").Value = -3182

$ws.Range("H131").Value = 795.80206
$ws.Range("I131").Value = 390
$ws.Range("J131").Value = 842.98834
$ws.Range("K131").Value = 1170
$ws.Range("L131").Value = 2528.96502
$ws.Range("M131").Value = 3870
$ws.Range("N131").Value = -12608.96502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 59800
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 59800
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 59800
$ws.Range("N74").Value = -61672

$ws.Range("H77").Value = 59800
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 59800
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 179400
$ws.Range("N77").Value = -188760

$ws.Range("H97").Value = 1865
$ws.Range("I97").Value = 2126
$ws.Range("J97").Value = 1212.5
$ws.Range("K97").Value = 2126
$ws.Range("L97").Value = 1212.5
$ws.Range("M97").Value = -1630
$ws.Range("N97").Value = -2204.5

$ws.Range("H126").Value = 2165.6428
$ws.Range("I126").Value = 2264.75
$ws.Range("J126").Value = 2033.5
$ws.Range("K126").Value = 6794.25
$ws.Range("L126").Value = 6100.5
$ws.Range("M126").Value = -4324.25
$ws.Range("N126").Value = -11040.5

$ws.Range("H132").Value = 2093
$ws.Range("I132").Value = 1608.0769
$ws.Range("J132").Value = 3143.6667
$ws.Range("K132").Value = 4824.2307
$ws.Range("L132").Value = 9431.000100000001
$ws.Range("M132").Value = -2294.2307
$ws.Range("N132").Value = -14491.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11433.429
$ws.Range("I132").Value = 13234.857
$ws.Range("J132").Value = 8731.286
$ws.Range("K132").Value = 39704.571
$ws.Range("L132").Value = 26193.858
$ws.Range("M132").Value = -37174.571
$ws.Range("N132").Value = -31253.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 600
$ws.Range("I4").Value = 600
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -487
$ws.Range("N4").ClearContents()

$ws.Range("H132").Value = 3218.8572
$ws.Range("I132").Value = 3022.5
$ws.Range("J132").Value = 3415.2144
$ws.Range("K132").Value = 9067.5
$ws.Range("L132").Value = 10245.6432
$ws.Range("M132").Value = -6537.5
$ws.Range("N132").Value = -15305.6432
